$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Oikeudet")

# Insert a new row at position 46 (shifts rows 46.. down by one)
$ws.Rows.Item(46).Insert()

# Populate the new row's cells
$ws.Cells.Item(46, 1).Value = "Urakat"
$ws.Cells.Item(46, 2).Value = "Kanavat / Kanavakohteet"
$ws.Cells.Item(46, 4).Value = "R*"
$ws.Cells.Item(46, 25).Value = "Tällä oikeudella voi hakea urakan kanavakohteet (tarvitaan mm. lomakkeissa)"
